# Add 2022-Q1 sheet with fund holdings data, and update the 总计 (Total) summary sheet.

$wb = $excel.ActiveWorkbook

# --- Fund holdings data for the new 2022-Q1 sheet ---
# Columns: 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$fundData = @(
    @("118001","易方达亚洲精选股票(QDII)","46.85","93.54","7.70","3.6074",7),
    @("008283","易方达金融行业股票","34.73","91.49","6.67","2.3165",7),
    @("011152","富兰克林国海兴海回报混合","17.18","84.91","5.34","0.9174",3),
    @("008515","富兰克林国海基本面优选混合","13.82","89.34","5.89","0.8140",5),
    @("011468","富兰克林国海竞争优势三年持有期混合型证券投资基金A","13.11","85.53","4.27","0.5598",6),
    @("870009","广发资管平衡精选一年持有混合A","11.34","94.29","4.73","0.5364",7),
    @("501025","鹏华港股通中证香港银行投资指数（LOF）A","9.81","94.47","3.40","0.3335",10),
    @("011913","华夏永泓一年持有混合A","24.65","37.51","1.04","0.2564",8),
    @("010365","鹏华港股通中证香港银行投资指数（LOF）C","6.07","94.47","3.40","0.2064",10),
    @("011914","华夏永泓一年持有混合C","12.39","37.51","1.04","0.1289",8),
    @("013009","万家港股通精选混合A","2.97","81.62","3.98","0.1182",8),
    @("160125","南方香港优选股票QDII-LOF","2.46","91.14","3.73","0.0918",3),
    @("012170","华夏永顺一年持有期混合型证券投资基金A","10.95","29.59","0.82","0.0898",8),
    @("161229","国投瑞银中国价值发现股票QDII-LOF","1.47","92.83","5.22","0.0767",6),
    @("872019","广发资管平衡精选一年持有混合C","1.54","94.29","4.73","0.0728",7),
    @("006809","泰康港股通中证香港银行投资指数A","1.99","94.73","3.40","0.0677",10),
    @("007354","创金合信港股通量化股票A","3.84","91.20","1.60","0.0614",10),
    @("241001","华宝海外中国混合(QDII)","0.83","86.89","6.27","0.0520",5),
    @("009017","银华港股通精选股票","0.91","86.12","4.77","0.0434",7),
    @("013010","万家港股通精选混合C","0.86","81.62","3.98","0.0342",8),
    @("006810","泰康港股通中证香港银行投资指数C","0.90","94.73","3.40","0.0306",10),
    @("011469","富兰克林国海竞争优势三年持有期混合型证券投资基金C","0.70","85.53","4.27","0.0299",6),
    @("006781","汇丰晋信港股通精选股票","0.67","90.36","3.95","0.0265",4),
    @("011647","博时港股通红利精选混合A","0.13","92.10","6.72","0.0087",4),
    @("001942","前海开源沪港深汇鑫灵活配置混合A","0.10","90.39","7.93","0.0079",4),
    @("001943","前海开源沪港深汇鑫灵活配置混合C","0.08","90.39","7.93","0.0063",4),
    @("012171","华夏永顺一年持有期混合型证券投资基金C","0.53","29.59","0.82","0.0043",8),
    @("007357","创金合信港股通量化股票C","0.26","91.20","1.60","0.0042",10),
    @("011648","博时港股通红利精选混合C","0.02","92.10","6.72","0.0013",4)
)

# --- 1. Insert the new "2022-Q1" worksheet right before the "总计" sheet ---
# NOTE: a Worksheet reference captured before a Worksheets.Add() call tracks
# the *position*, not the original sheet identity, once sheets shift - so the
# "总计" handle must be re-fetched by name after the insert.
$templateSheet = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

$totalSheet = $wb.Worksheets.Item("总计")

# --- 2. Write header row (B1:H1) ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- 3. Write data rows (A2:H30) ---
$rowCount = $fundData.Count
$lastRow = 1 + $rowCount

$textRange = $newSheet.Range("B2:G" + $lastRow)
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $item = $fundData[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = $item[0]
    $newSheet.Cells.Item($r, 3).Value = $item[1]
    $newSheet.Cells.Item($r, 4).Value = $item[2]
    $newSheet.Cells.Item($r, 5).Value = $item[3]
    $newSheet.Cells.Item($r, 6).Value = $item[4]
    $newSheet.Cells.Item($r, 7).Value = $item[5]
    $newSheet.Cells.Item($r, 8).Value = $item[6]
}

# Remove the auto-applied text NumberFormat so cells carry no explicit style,
# matching the plain inline-string cells produced by the source pipeline.
$textRange.Style = "Normal"

# --- 4. Re-apply the header/index styling (bold, centered, bordered) used
#        throughout the workbook, copied from the template sheet ---
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A" + $lastRow).PasteSpecial(-4122)

# --- 5. Rebuild the "总计" (Total) summary sheet with the new 2022-Q1 row ---
$totalSheet.Cells.Clear()

$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$totalData = @(
    @("2022-Q1", 29, 10.5),
    @("2021-Q4", 19, 9.98),
    @("2021-Q3", 17, 9.970000000000001),
    @("2021-Q2", 16, 9.08),
    @("2021-Q1", 22, 7.07),
    @("2020-Q4", 12, 2.66)
)

for ($i = 0; $i -lt $totalData.Count; $i++) {
    $r = $i + 2
    $item = $totalData[$i]
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $item[0]
    $totalSheet.Cells.Item($r, 3).Value = $item[1]
    $totalSheet.Cells.Item($r, 4).Value = $item[2]
}

$templateSheet.Range("B1:D1").Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122)

$templateSheet.Range("A2").Copy()
$totalSheet.Range("A2:A" + (1 + $totalData.Count)).PasteSpecial(-4122)

$newSheet.Range("A1").Select()
$totalSheet.Range("A1").Select()
